$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the data (columns F..V) between the following row pairs, leaving
#    columns A (index) and B-E (league/date metadata, identical anyway)
#    untouched. This reorders mis-ordered match rows back into the
#    canonical (chronological) order recorded by the scraper.
# ---------------------------------------------------------------------------
$pairs = @(
    @(10,11),
    @(15,16),
    @(29,30),
    @(71,72),
    @(87,88),
    @(89,90),
    @(91,92)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = "F" + $r1 + ":V" + $r1
    $range2 = "F" + $r2 + ":V" + $r2
    $v1 = $ws.Range($range1).Value2
    $v2 = $ws.Range($range2).Value2
    $ws.Range($range1).Value2 = $v2
    $ws.Range($range2).Value2 = $v1
}

# ---------------------------------------------------------------------------
# 2) Append three new match rows (114-116) scraped by the 11/11/2023 run.
# ---------------------------------------------------------------------------
$newRows = @(
    @(114, 113, "turkey", "super-lig", "2023-2024", 45241.47916666666, "Kasimpasa", 3, "Kayserispor", 4, 2.31, "06/11/2023 18:12", 2.33, "11/11/2023 11:21", 3.64, "06/11/2023 18:12", 3.71, "11/11/2023 11:18", 3.02, "06/11/2023 18:12", 3.04, "11/11/2023 11:21", "https://www.betexplorer.com/football/turkey/super-lig/kasimpasa-kayserispor/E72uxNBI/"),
    @(115, 114, "turkey", "super-lig", "2023-2024", 45241.58333333334, "Alanyaspor", 0, "Gaziantep", 3, 2.07, "06/11/2023 18:12", 2.34, "11/11/2023 13:59", 3.62, "06/11/2023 18:12", 3.29, "11/11/2023 13:55", 3.67, "06/11/2023 18:12", 3.37, "11/11/2023 13:59", "https://www.betexplorer.com/football/turkey/super-lig/alanyaspor-gaziantep/hb1qysRO/"),
    @(116, 115, "turkey", "super-lig", "2023-2024", 45241.58333333334, "Karagumruk", 3, "Sivasspor", 0, 2.41, "06/11/2023 04:12", 2.22, "11/11/2023 13:53", 3.41, "06/11/2023 04:12", 3.42, "11/11/2023 13:53", 3.1, "06/11/2023 04:12", 3.51, "11/11/2023 13:53", "https://www.betexplorer.com/football/turkey/super-lig/f-karagumruk-sivasspor/WOWgepCt/")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
    $ws.Cells.Item($r, 14).Value = $row[14]
    $ws.Cells.Item($r, 15).Value = $row[15]
    $ws.Cells.Item($r, 16).Value = $row[16]
    $ws.Cells.Item($r, 17).Value = $row[17]
    $ws.Cells.Item($r, 18).Value = $row[18]
    $ws.Cells.Item($r, 19).Value = $row[19]
    $ws.Cells.Item($r, 20).Value = $row[20]
    $ws.Cells.Item($r, 21).Value = $row[21]
    $ws.Cells.Item($r, 22).Value = $row[22]

    # Match the formatting used by every other data row: column A (index)
    # carries the bold/centered/bordered header-like style, column E (match
    # date) carries the custom date-time number format.
    $ws.Range("A113").Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)
    $ws.Range("E113").Copy()
    $ws.Range("E" + $r).PasteSpecial(-4122)
}
